$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for rule R10 (cell E8) as committed via "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the saved selection/active cell in the sheet view
$ws.Range("E8").Select()
